$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$ws.Range("A1").Font.Name = "Calibri"
